{"js": "// Remove the trailing \"Ver no Jupiter ...\" / copyright footer block (and the\n// blank paragraph that separated it from the bibliography) that the Jekyll\n// site build no longer emits at the end of the document.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph (\"...uma introdu\u00e7\u00e3o \u00e0 \u00e1lgebra linear...\") so we\n// only touch the footer block that immediately follows it, not any other\n// blank paragraph elsewhere in the document.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"uma introdu\u00e7\u00e3o \u00e0 \u00e1lgebra linear\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  // The paragraphs right after the anchor are expected to be:\n  //   (empty), \"Ver no Jupiter ...\", \"\u00a9 2020 ... Creative Commons Attribution\"\n  // Delete exactly those three paragraphs.\n  const toDelete = [];\n  let idx = anchorIndex + 1;\n  if (idx < items.length && items[idx].text === \"\") {\n    toDelete.push(items[idx]);\n    idx++;\n  }\n  if (idx < items.length && items[idx].text.indexOf(\"Ver no Jupiter\") !== -1) {\n    toDelete.push(items[idx]);\n    idx++;\n  }\n  if (idx < items.length && items[idx].text.indexOf(\"Powered by Jekyll\") !== -1) {\n    toDelete.push(items[idx]);\n    idx++;\n  }\n\n  for (const p of toDelete) {\n    p.delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / copyright footer block (and the\n# blank paragraph that separated it from the bibliography) that the Jekyll\n# site build no longer emits at the end of the document.\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph (\"...uma introdu\u00e7\u00e3o \u00e0 \u00e1lgebra linear... Thomson,\n# 2007.\") so we only touch the footer block that immediately follows it, not\n# any other blank paragraph elsewhere in the document.\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n  $t = $d.Paragraphs.Item($i).Range.Text\n  if ($t -like \"*Thomson, 2007*\") {\n    $anchorIndex = $i\n    break\n  }\n}\n\nif ($anchorIndex -gt 0) {\n  # The paragraphs right after the anchor are expected to be:\n  #   (empty), \"Ver no Jupiter ...\", \"\u00a9 2020 ... Creative Commons Attribution\"\n  # Collect exactly those, then delete back-to-front so earlier indices stay\n  # valid while later ones are removed.\n  $targets = @()\n  $idx = $anchorIndex + 1\n\n  if ($idx -le $d.Paragraphs.Count) {\n    $t = $d.Paragraphs.Item($idx).Range.Text\n    if ($t.Trim() -eq \"\") {\n      $targets += $idx\n      $idx = $idx + 1\n    }\n  }\n  if ($idx -le $d.Paragraphs.Count) {\n    $t = $d.Paragraphs.Item($idx).Range.Text\n    if ($t -like \"*Ver no Jupiter*\") {\n      $targets += $idx\n      $idx = $idx + 1\n    }\n  }\n  if ($idx -le $d.Paragraphs.Count) {\n    $t = $d.Paragraphs.Item($idx).Range.Text\n    if ($t -like \"*Powered by Jekyll*\") {\n      $targets += $idx\n      $idx = $idx + 1\n    }\n  }\n\n  for ($j = $targets.Count - 1; $j -ge 0; $j--) {\n    $d.Paragraphs.Item($targets[$j]).Range.Delete()\n  }\n}\n"}
